$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused rows (6 through 11) so the sheet shrinks to A1:B5
$ws.Rows("6:11").Delete()

# Header row
$ws.Range("A1").Value = "keyword"
$ws.Range("B1").Value = "answer"

# Keywords for each row
$ws.Range("A2").Value = "wifi"
$ws.Range("A3").Value = "mytv"
$ws.Range("A4").Value = "mesh"
$ws.Range("A5").Value = "sim"

# Answers for each row
$ws.Range("B2").Value = "Chúng tôi hỗ trợ sửa chữa và lắp đặt Wifi. Gọi 18001091 để biết thêm chi tiết."
$ws.Range("B3").Value = "Dịch vụ MyTV hiện có nhiều gói hấp dẫn. Bạn muốn nâng cấp hay sửa chữa?"
$ws.Range("B4").Value = "Hệ thống Mesh giúp mở rộng sóng wifi. Bạn cần tư vấn thêm không?"
$ws.Range("B5").Value = "Chúng tôi cung cấp SIM Vinaphone chính hãng. Bạn muốn mua gói nào?"
